# Generate Report for Handoff
# Adds two new handed-off files (40c131f6-...  and 53fa06c0-...) as new
# rows (4 and 5) on all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# Hyperlink-style colour used by the workbook's existing "HyperLink" cell
# style (RGB 64,95,ED == OLE colour 15570276) plus a single underline, so
# newly-hyperlinked cells look like the pre-existing ones.
$hyperlinkColor = 15570276
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $hyperlinkColor
}

function Style-AsDate($rng) {
    $rng.NumberFormat = $dateFormat
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 4 - 40c131f6-07a7-4fde-8562-5a5241c1fb5b
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ba102612fc082eb1f84afb3f944e0b0e789f03da/e2e/40c131f6-07a7-4fde-8562-5a5241c1fb5b.md", "", "", "40c131f6-07a7-4fde-8562-5a5241c1fb5b.md") | Out-Null
Style-AsHyperlink $wsOverview.Range("A4")
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-03-23 08:42:52"
Style-AsDate $wsOverview.Range("D4")

# Row 5 - 53fa06c0-b04a-41e0-9992-9ecea00039f5
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/aff86616d7d4e538e687e7848667853528cc46da/e2e/53fa06c0-b04a-41e0-9992-9ecea00039f5.md", "", "", "53fa06c0-b04a-41e0-9992-9ecea00039f5.md") | Out-Null
Style-AsHyperlink $wsOverview.Range("A5")
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-03-23 08:42:52"
Style-AsDate $wsOverview.Range("D5")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 4 - 40c131f6-07a7-4fde-8562-5a5241c1fb5b
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ba102612fc082eb1f84afb3f944e0b0e789f03da/e2e/40c131f6-07a7-4fde-8562-5a5241c1fb5b.md", "", "", "40c131f6-07a7-4fde-8562-5a5241c1fb5b.md") | Out-Null
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ba102612fc082eb1f84afb3f944e0b0e789f03da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/40c131f6-07a7-4fde-8562-5a5241c1fb5b.ba102612fc082eb1f84afb3f944e0b0e789f03da.zh-cn.xlf", "", "", "40c131f6-07a7-4fde-8562-5a5241c1fb5b.ba102612fc082eb1f84afb3f944e0b0e789f03da.zh-cn.xlf") | Out-Null
$wsZh.Range("E4").Value = "2016-03-23 08:42:48"
Style-AsDate $wsZh.Range("E4")
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
Style-AsDate $wsZh.Range("H4")
$wsZh.Range("J4").Value = "Include"

# Row 5 - 53fa06c0-b04a-41e0-9992-9ecea00039f5
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/aff86616d7d4e538e687e7848667853528cc46da/e2e/53fa06c0-b04a-41e0-9992-9ecea00039f5.md", "", "", "53fa06c0-b04a-41e0-9992-9ecea00039f5.md") | Out-Null
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aff86616d7d4e538e687e7848667853528cc46da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/53fa06c0-b04a-41e0-9992-9ecea00039f5.aff86616d7d4e538e687e7848667853528cc46da.zh-cn.xlf", "", "", "53fa06c0-b04a-41e0-9992-9ecea00039f5.aff86616d7d4e538e687e7848667853528cc46da.zh-cn.xlf") | Out-Null
$wsZh.Range("E5").Value = "2016-03-23 08:42:48"
Style-AsDate $wsZh.Range("E5")
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
Style-AsDate $wsZh.Range("H5")
$wsZh.Range("J5").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 4 - 40c131f6-07a7-4fde-8562-5a5241c1fb5b
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ba102612fc082eb1f84afb3f944e0b0e789f03da/e2e/40c131f6-07a7-4fde-8562-5a5241c1fb5b.md", "", "", "40c131f6-07a7-4fde-8562-5a5241c1fb5b.md") | Out-Null
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ba102612fc082eb1f84afb3f944e0b0e789f03da/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/40c131f6-07a7-4fde-8562-5a5241c1fb5b.ba102612fc082eb1f84afb3f944e0b0e789f03da.de-de.xlf", "", "", "40c131f6-07a7-4fde-8562-5a5241c1fb5b.ba102612fc082eb1f84afb3f944e0b0e789f03da.de-de.xlf") | Out-Null
$wsDe.Range("E4").Value = "2016-03-23 08:42:52"
Style-AsDate $wsDe.Range("E4")
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
Style-AsDate $wsDe.Range("H4")
$wsDe.Range("J4").Value = "Include"

# Row 5 - 53fa06c0-b04a-41e0-9992-9ecea00039f5
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/aff86616d7d4e538e687e7848667853528cc46da/e2e/53fa06c0-b04a-41e0-9992-9ecea00039f5.md", "", "", "53fa06c0-b04a-41e0-9992-9ecea00039f5.md") | Out-Null
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aff86616d7d4e538e687e7848667853528cc46da/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/53fa06c0-b04a-41e0-9992-9ecea00039f5.aff86616d7d4e538e687e7848667853528cc46da.de-de.xlf", "", "", "53fa06c0-b04a-41e0-9992-9ecea00039f5.aff86616d7d4e538e687e7848667853528cc46da.de-de.xlf") | Out-Null
$wsDe.Range("E5").Value = "2016-03-23 08:42:52"
Style-AsDate $wsDe.Range("E5")
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
Style-AsDate $wsDe.Range("H5")
$wsDe.Range("J5").Value = "Include"

Write-Host "Done: added rows 4-5 to Overview, zh-cn, de-de sheets."
